$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ideal-format")

# Replace the static N values with formulas that sum the component rows
$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"
$ws.Range("C3").Formula = "=C16+C17+C18+C19+C20+C21+C22"

# Update the active selection to reflect the author's new cursor position
$ws.Range("C4").Select()
